# Add PMID column to "studies" sheet and notes column to "counts" sheet
$wb = $excel.ActiveWorkbook

$wsStudies = $wb.Worksheets.Item("studies")
$wsCounts  = $wb.Worksheets.Item("counts")

# --- studies sheet: add new "PMID" header in column H ---
$wsStudies.Activate()
$wsStudies.Range("H1").Value = "PMID"
$wsStudies.Range("H2").Select()

# --- counts sheet: add new "notes" header in column F ---
$wsCounts.Activate()
$wsCounts.Range("F1").Value = "notes"
$wsCounts.Range("F2").Select()

# leave "counts" sheet as the active tab, matching the target workbook state
$wsCounts.Activate()
